$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Swap the F:V (match detail) content between duplicated-timestamp row pairs ---
# (columns A-E -- Indice/pais/torneio/temporada/data_partida -- stay put per row)
$rowA = $ws.Range("F3:V3").Value2
$rowB = $ws.Range("F4:V4").Value2
$ws.Range("F3:V3").Value2 = $rowB
$ws.Range("F4:V4").Value2 = $rowA

$rowA = $ws.Range("F18:V18").Value2
$rowB = $ws.Range("F19:V19").Value2
$ws.Range("F18:V18").Value2 = $rowB
$ws.Range("F19:V19").Value2 = $rowA

$rowA = $ws.Range("F31:V31").Value2
$rowB = $ws.Range("F32:V32").Value2
$ws.Range("F31:V31").Value2 = $rowB
$ws.Range("F32:V32").Value2 = $rowA

$rowA = $ws.Range("F33:V33").Value2
$rowB = $ws.Range("F34:V34").Value2
$ws.Range("F33:V33").Value2 = $rowB
$ws.Range("F34:V34").Value2 = $rowA

$rowA = $ws.Range("F42:V42").Value2
$rowB = $ws.Range("F43:V43").Value2
$ws.Range("F42:V42").Value2 = $rowB
$ws.Range("F43:V43").Value2 = $rowA

$rowA = $ws.Range("F44:V44").Value2
$rowB = $ws.Range("F45:V45").Value2
$ws.Range("F44:V44").Value2 = $rowB
$ws.Range("F45:V45").Value2 = $rowA

$rowA = $ws.Range("F51:V51").Value2
$rowB = $ws.Range("F52:V52").Value2
$ws.Range("F51:V51").Value2 = $rowB
$ws.Range("F52:V52").Value2 = $rowA

$rowA = $ws.Range("F55:V55").Value2
$rowB = $ws.Range("F56:V56").Value2
$ws.Range("F55:V55").Value2 = $rowB
$ws.Range("F56:V56").Value2 = $rowA

$rowA = $ws.Range("F68:V68").Value2
$rowB = $ws.Range("F69:V69").Value2
$ws.Range("F68:V68").Value2 = $rowB
$ws.Range("F69:V69").Value2 = $rowA

$rowA = $ws.Range("F82:V82").Value2
$rowB = $ws.Range("F83:V83").Value2
$ws.Range("F82:V82").Value2 = $rowB
$ws.Range("F83:V83").Value2 = $rowA

# --- Append 8 new match rows (122-129), matching formatting of the last existing row ---
$ws.Range("A121:V121").Copy()
$ws.Range("A122:V129").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# row 122
$ws.Range("A122").Value2 = 121
$ws.Range("B122").Value2 = "turkey"
$ws.Range("C122").Value2 = "super-lig"
$ws.Range("D122").Value2 = "2023-2024"
$ws.Range("E122").Value2 = 45255.47916666666
$fv = New-Object 'object[,]' 1,17
$fv[0,0] = "Istanbulspor AS"
$fv[0,1] = 2
$fv[0,2] = "Hatayspor"
$fv[0,3] = 1
$fv[0,4] = 2.67
$fv[0,5] = "12/11/2023 15:42"
$fv[0,6] = 3.12
$fv[0,7] = "25/11/2023 11:29"
$fv[0,8] = 3.46
$fv[0,9] = "12/11/2023 15:42"
$fv[0,10] = 3.36
$fv[0,11] = "25/11/2023 11:27"
$fv[0,12] = 2.67
$fv[0,13] = "12/11/2023 15:42"
$fv[0,14] = 2.44
$fv[0,15] = "25/11/2023 11:29"
$fv[0,16] = "https://www.betexplorer.com/football/turkey/super-lig/istanbulspor-as-hatayspor/YVyIkttI/"
$ws.Range("F122:V122").Value2 = $fv

# row 123
$ws.Range("A123").Value2 = 122
$ws.Range("B123").Value2 = "turkey"
$ws.Range("C123").Value2 = "super-lig"
$ws.Range("D123").Value2 = "2023-2024"
$ws.Range("E123").Value2 = 45255.58333333334
$fv = New-Object 'object[,]' 1,17
$fv[0,0] = "Gaziantep"
$fv[0,1] = 0
$fv[0,2] = "Ankaragucu"
$fv[0,3] = 1
$fv[0,4] = 2.14
$fv[0,5] = "12/11/2023 15:42"
$fv[0,6] = 2.09
$fv[0,7] = "25/11/2023 13:56"
$fv[0,8] = 3.57
$fv[0,9] = "12/11/2023 15:42"
$fv[0,10] = 3.4
$fv[0,11] = "25/11/2023 13:56"
$fv[0,12] = 3.51
$fv[0,13] = "12/11/2023 15:42"
$fv[0,14] = 3.89
$fv[0,15] = "25/11/2023 13:59"
$fv[0,16] = "https://www.betexplorer.com/football/turkey/super-lig/gaziantep-ankaragucu/6J6W83tt/"
$ws.Range("F123:V123").Value2 = $fv

# row 124
$ws.Range("A124").Value2 = 123
$ws.Range("B124").Value2 = "turkey"
$ws.Range("C124").Value2 = "super-lig"
$ws.Range("D124").Value2 = "2023-2024"
$ws.Range("E124").Value2 = 45255.58333333334
$fv = New-Object 'object[,]' 1,17
$fv[0,0] = "Konyaspor"
$fv[0,1] = 2
$fv[0,2] = "Kasimpasa"
$fv[0,3] = 0
$fv[0,4] = 1.99
$fv[0,5] = "12/11/2023 15:42"
$fv[0,6] = 2.08
$fv[0,7] = "25/11/2023 13:56"
$fv[0,8] = 3.77
$fv[0,9] = "12/11/2023 15:42"
$fv[0,10] = 3.68
$fv[0,11] = "25/11/2023 13:56"
$fv[0,12] = 3.78
$fv[0,13] = "12/11/2023 15:42"
$fv[0,14] = 3.62
$fv[0,15] = "25/11/2023 13:56"
$fv[0,16] = "https://www.betexplorer.com/football/turkey/super-lig/konyaspor-kasimpasa/nZ1i4ueB/"
$ws.Range("F124:V124").Value2 = $fv

# row 125
$ws.Range("A125").Value2 = 124
$ws.Range("B125").Value2 = "turkey"
$ws.Range("C125").Value2 = "super-lig"
$ws.Range("D125").Value2 = "2023-2024"
$ws.Range("E125").Value2 = 45255.70833333334
$fv = New-Object 'object[,]' 1,17
$fv[0,0] = "Galatasaray"
$fv[0,1] = 4
$fv[0,2] = "Alanyaspor"
$fv[0,3] = 0
$fv[0,4] = 1.19
$fv[0,5] = "12/11/2023 15:42"
$fv[0,6] = 1.21
$fv[0,7] = "25/11/2023 16:59"
$fv[0,8] = 8.109999999999999
$fv[0,9] = "12/11/2023 15:42"
$fv[0,10] = 7.57
$fv[0,11] = "25/11/2023 16:59"
$fv[0,12] = 14.43
$fv[0,13] = "12/11/2023 15:42"
$fv[0,14] = 13.33
$fv[0,15] = "25/11/2023 16:59"
$fv[0,16] = "https://www.betexplorer.com/football/turkey/super-lig/galatasaray-alanyaspor/GfRMl0eO/"
$ws.Range("F125:V125").Value2 = $fv

# row 126
$ws.Range("A126").Value2 = 125
$ws.Range("B126").Value2 = "turkey"
$ws.Range("C126").Value2 = "super-lig"
$ws.Range("D126").Value2 = "2023-2024"
$ws.Range("E126").Value2 = 45256.47916666666
$fv = New-Object 'object[,]' 1,17
$fv[0,0] = "Kayserispor"
$fv[0,1] = 1
$fv[0,2] = "Adana Demirspor"
$fv[0,3] = 1
$fv[0,4] = 2.61
$fv[0,5] = "12/11/2023 17:12"
$fv[0,6] = 2.29
$fv[0,7] = "26/11/2023 11:26"
$fv[0,8] = 3.69
$fv[0,9] = "12/11/2023 17:12"
$fv[0,10] = 3.88
$fv[0,11] = "26/11/2023 11:20"
$fv[0,12] = 2.61
$fv[0,13] = "12/11/2023 17:12"
$fv[0,14] = 3
$fv[0,15] = "26/11/2023 11:26"
$fv[0,16] = "https://www.betexplorer.com/football/turkey/super-lig/kayserispor-adanademirspor/fuwEjMQB/"
$ws.Range("F126:V126").Value2 = $fv

# row 127
$ws.Range("A127").Value2 = 126
$ws.Range("B127").Value2 = "turkey"
$ws.Range("C127").Value2 = "super-lig"
$ws.Range("D127").Value2 = "2023-2024"
$ws.Range("E127").Value2 = 45256.58333333334
$fv = New-Object 'object[,]' 1,17
$fv[0,0] = "Antalyaspor"
$fv[0,1] = 0
$fv[0,2] = "Rizespor"
$fv[0,3] = 0
$fv[0,4] = 1.67
$fv[0,5] = "12/11/2023 15:42"
$fv[0,6] = 1.68
$fv[0,7] = "26/11/2023 13:57"
$fv[0,8] = 4.06
$fv[0,9] = "12/11/2023 15:42"
$fv[0,10] = 3.88
$fv[0,11] = "26/11/2023 13:55"
$fv[0,12] = 5.07
$fv[0,13] = "12/11/2023 15:42"
$fv[0,14] = 5.55
$fv[0,15] = "26/11/2023 13:55"
$fv[0,16] = "https://www.betexplorer.com/football/turkey/super-lig/antalyaspor-rizespor/2o0u7sBh/"
$ws.Range("F127:V127").Value2 = $fv

# row 128
$ws.Range("A128").Value2 = 127
$ws.Range("B128").Value2 = "turkey"
$ws.Range("C128").Value2 = "super-lig"
$ws.Range("D128").Value2 = "2023-2024"
$ws.Range("E128").Value2 = 45256.58333333334
$fv = New-Object 'object[,]' 1,17
$fv[0,0] = "Samsunspor"
$fv[0,1] = 1
$fv[0,2] = "Besiktas"
$fv[0,3] = 2
$fv[0,4] = 3.24
$fv[0,5] = "12/11/2023 15:42"
$fv[0,6] = 3.31
$fv[0,7] = "26/11/2023 13:57"
$fv[0,8] = 3.56
$fv[0,9] = "12/11/2023 15:42"
$fv[0,10] = 3.5
$fv[0,11] = "26/11/2023 13:57"
$fv[0,12] = 2.23
$fv[0,13] = "12/11/2023 15:42"
$fv[0,14] = 2.27
$fv[0,15] = "26/11/2023 13:57"
$fv[0,16] = "https://www.betexplorer.com/football/turkey/super-lig/samsunspor-besiktas/UcVQmKAU/"
$ws.Range("F128:V128").Value2 = $fv

# row 129
$ws.Range("A129").Value2 = 128
$ws.Range("B129").Value2 = "turkey"
$ws.Range("C129").Value2 = "super-lig"
$ws.Range("D129").Value2 = "2023-2024"
$ws.Range("E129").Value2 = 45256.70833333334
$fv = New-Object 'object[,]' 1,17
$fv[0,0] = "Fenerbahce"
$fv[0,1] = 2
$fv[0,2] = "Karagumruk"
$fv[0,3] = 1
$fv[0,4] = 1.25
$fv[0,5] = "13/11/2023 02:42"
$fv[0,6] = 1.25
$fv[0,7] = "26/11/2023 16:58"
$fv[0,8] = 6.63
$fv[0,9] = "13/11/2023 02:42"
$fv[0,10] = 6.7
$fv[0,11] = "26/11/2023 16:58"
$fv[0,12] = 10.71
$fv[0,13] = "13/11/2023 02:42"
$fv[0,14] = 11.72
$fv[0,15] = "26/11/2023 16:58"
$fv[0,16] = "https://www.betexplorer.com/football/turkey/super-lig/fenerbahce-f-karagumruk/txaq61Qb/"
$ws.Range("F129:V129").Value2 = $fv

$ws.Range("A1").Select() | Out-Null
